$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.44367991736487372
$ws.Range("G2").Value = 3101.4848772818259
$ws.Range("L2").Value = 0.1065
$ws.Range("S2").Value = 0.55249752902628602

$ws.Range("K6").Select()
